$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value2 = 1079.4  # H41: 1649.5 -> 1079.4
$ws.Cells.Item(41, 9).Value2 = 99  # I41: 2299.5 -> 99
$ws.Cells.Item(41, 11).Value2 = 99  # K41: 2299.5 -> 99
$ws.Cells.Item(41, 13).Value2 = 341  # M41: -1859.5 -> 341
$ws.Cells.Item(62, 8).Value2 = 7548.4614  # H62: 8052.6 -> 7548.4614
$ws.Cells.Item(62, 9).Value2 = 6566  # I62: 6949.5 -> 6566
$ws.Cells.Item(62, 10).Value2 = 8390.571  # J62: 8788 -> 8390.571
$ws.Cells.Item(62, 11).Value2 = 6566  # K62: 6949.5 -> 6566
$ws.Cells.Item(62, 12).Value2 = 8390.571  # L62: 8788 -> 8390.571
$ws.Cells.Item(62, 13).Value2 = -5942  # M62: -6325.5 -> -5942
$ws.Cells.Item(62, 14).Value2 = -9638.571  # N62: -10036 -> -9638.571
$ws.Cells.Item(65, 8).Value2 = 7548.4614  # H65: 8052.6 -> 7548.4614
$ws.Cells.Item(65, 9).Value2 = 6566  # I65: 6949.5 -> 6566
$ws.Cells.Item(65, 10).Value2 = 8390.571  # J65: 8788 -> 8390.571
$ws.Cells.Item(65, 11).Value2 = 32830  # K65: 34747.5 -> 32830
$ws.Cells.Item(65, 12).Value2 = 41952.855  # L65: 43940 -> 41952.855
$ws.Cells.Item(65, 13).Value2 = -29710  # M65: -31627.5 -> -29710
$ws.Cells.Item(65, 14).Value2 = -48192.855  # N65: -50180 -> -48192.855
$ws.Cells.Item(94, 8).Value2 = 9007.154  # H94: 8826.182000000001 -> 9007.154
$ws.Cells.Item(94, 9).Value2 = 9007.154  # I94: 8826.182000000001 -> 9007.154
$ws.Cells.Item(94, 11).Value2 = 9007.154  # K94: 8826.182000000001 -> 9007.154
$ws.Cells.Item(94, 13).Value2 = -8556.154  # M94: -8375.182000000001 -> -8556.154
$ws.Cells.Item(106, 8).Value2 = 924.6667  # H106: 1137.5 -> 924.6667
$ws.Cells.Item(106, 9).Value2 = 924.6667  # I106: 1137.5 -> 924.6667
$ws.Cells.Item(106, 11).Value2 = 924.6667  # K106: 1137.5 -> 924.6667
$ws.Cells.Item(106, 13).Value2 = -293.6667  # M106: -506.5 -> -293.6667
$ws.Cells.Item(107, 8).Value2 = 231.6923  # H107: 234.42308 -> 231.6923
$ws.Cells.Item(107, 9).Value2 = 94.57143000000001  # I107: 97.95238000000001 -> 94.57143000000001
$ws.Cells.Item(107, 11).Value2 = 94.57143000000001  # K107: 97.95238000000001 -> 94.57143000000001
$ws.Cells.Item(107, 13).Value2 = 1825.42857  # M107: 1822.04762 -> 1825.42857
$ws.Cells.Item(125, 8).Value2 = 2828.3333  # H125: 2992.5 -> 2828.3333
$ws.Cells.Item(125, 9).Value2 = 2828.3333  # I125: 2992.5 -> 2828.3333
$ws.Cells.Item(125, 11).Value2 = 25454.9997  # K125: 26932.5 -> 25454.9997
$ws.Cells.Item(125, 13).Value2 = -22994.9997  # M125: -24472.5 -> -22994.9997

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value2 = 3851  # H45: 4082 -> 3851
$ws.Cells.Item(45, 9).Value2 = 3079.1  # I45: 3261.5 -> 3079.1
$ws.Cells.Item(45, 11).Value2 = 3079.1  # K45: 3261.5 -> 3079.1
$ws.Cells.Item(45, 13).Value2 = -2702.1  # M45: -2884.5 -> -2702.1
$ws.Cells.Item(61, 8).Value2 = 4479.2  # H61: 4316 -> 4479.2
$ws.Cells.Item(61, 10).Value2 = 0  # J61: 3500 -> 0
$ws.Cells.Item(61, 12).Value2 = 0  # L61: 3500 -> 0
$ws.Cells.Item(61, 14).ClearContents()  # N61: -3924 -> (removed)
$ws.Cells.Item(132, 8).Value2 = 9992.5  # H132: 7654 -> 9992.5
$ws.Cells.Item(132, 9).Value2 = 9992.5  # I132: 7654 -> 9992.5
$ws.Cells.Item(132, 11).Value2 = 29977.5  # K132: 22962 -> 29977.5
$ws.Cells.Item(132, 13).Value2 = -27447.5  # M132: -20432 -> -27447.5
$ws.Cells.Item(136, 8).Value2 = 4479.2  # H136: 4316 -> 4479.2
$ws.Cells.Item(136, 10).Value2 = 0  # J136: 3500 -> 0
$ws.Cells.Item(136, 12).Value2 = 0  # L136: 10500 -> 0
$ws.Cells.Item(136, 14).ClearContents()  # N136: -15600 -> (removed)

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value2 = 23499.5  # H92: 23999.666 -> 23499.5
$ws.Cells.Item(92, 10).Value2 = 23499.5  # J92: 23999.666 -> 23499.5
$ws.Cells.Item(92, 12).Value2 = 23499.5  # L92: 23999.666 -> 23499.5
$ws.Cells.Item(92, 14).Value2 = -28491.5  # N92: -28991.666 -> -28491.5
$ws.Cells.Item(94, 8).Value2 = 1066.6666  # H94: 1100 -> 1066.6666
$ws.Cells.Item(105, 8).Value2 = 2254.8462  # H105: 2263.6667 -> 2254.8462
$ws.Cells.Item(105, 9).Value2 = 1874.25  # I105: 1856.75 -> 1874.25
$ws.Cells.Item(105, 10).Value2 = 2863.8  # J105: 3077.5 -> 2863.8
$ws.Cells.Item(105, 11).Value2 = 1874.25  # K105: 1856.75 -> 1874.25
$ws.Cells.Item(105, 12).Value2 = 2863.8  # L105: 3077.5 -> 2863.8
$ws.Cells.Item(105, 13).Value2 = -127.25  # M105: -109.75 -> -127.25
$ws.Cells.Item(105, 14).Value2 = -6357.8  # N105: -6571.5 -> -6357.8
$ws.Cells.Item(107, 8).Value2 = 3368.3547  # H107: 3374.1614 -> 3368.3547
$ws.Cells.Item(107, 9).Value2 = 1478.1364  # I107: 1486.3182 -> 1478.1364
$ws.Cells.Item(107, 11).Value2 = 1478.1364  # K107: 1486.3182 -> 1478.1364
$ws.Cells.Item(107, 13).Value2 = 441.8635999999999  # M107: 433.6818000000001 -> 441.8635999999999
$ws.Cells.Item(132, 8).Value2 = 160000  # H132: 0 -> 160000
$ws.Cells.Item(132, 10).Value2 = 160000  # J132: 0 -> 160000
$ws.Cells.Item(132, 12).Value2 = 160000  # L132: 0 -> 160000
$ws.Cells.Item(132, 14).Value2 = -170120  # N132: None -> -170120

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value2 = 9250  # H18: 9200 -> 9250
$ws.Cells.Item(18, 10).Value2 = 9250  # J18: 9200 -> 9250
$ws.Cells.Item(18, 12).Value2 = 9250  # L18: 9200 -> 9250
$ws.Cells.Item(18, 14).Value2 = -9710  # N18: -9660 -> -9710
$ws.Cells.Item(31, 8).Value2 = 5491.4287  # H31: 5395.7 -> 5491.4287
$ws.Cells.Item(31, 9).Value2 = 4522.5454  # I31: 4573.636 -> 4522.5454
$ws.Cells.Item(31, 10).Value2 = 9044  # J31: 7656.375 -> 9044
$ws.Cells.Item(31, 11).Value2 = 4522.5454  # K31: 4573.636 -> 4522.5454
$ws.Cells.Item(31, 12).Value2 = 9044  # L31: 7656.375 -> 9044
$ws.Cells.Item(31, 13).Value2 = -4227.5454  # M31: -4278.636 -> -4227.5454
$ws.Cells.Item(31, 14).Value2 = -9634  # N31: -8246.375 -> -9634
$ws.Cells.Item(34, 8).Value2 = 5491.4287  # H34: 5395.7 -> 5491.4287
$ws.Cells.Item(34, 9).Value2 = 4522.5454  # I34: 4573.636 -> 4522.5454
$ws.Cells.Item(34, 10).Value2 = 9044  # J34: 7656.375 -> 9044
$ws.Cells.Item(34, 11).Value2 = 4522.5454  # K34: 4573.636 -> 4522.5454
$ws.Cells.Item(34, 12).Value2 = 9044  # L34: 7656.375 -> 9044
$ws.Cells.Item(34, 13).Value2 = -4320.5454  # M34: -4371.636 -> -4320.5454
$ws.Cells.Item(34, 14).Value2 = -9448  # N34: -8060.375 -> -9448
$ws.Cells.Item(68, 8).Value2 = 91888.5  # H68: 92301.28999999999 -> 91888.5
$ws.Cells.Item(68, 10).Value2 = 91888.5  # J68: 92301.28999999999 -> 91888.5
$ws.Cells.Item(68, 12).Value2 = 91888.5  # L68: 92301.28999999999 -> 91888.5
$ws.Cells.Item(68, 14).Value2 = -93386.5  # N68: -93799.28999999999 -> -93386.5
$ws.Cells.Item(71, 8).Value2 = 91888.5  # H71: 92301.28999999999 -> 91888.5
$ws.Cells.Item(71, 10).Value2 = 91888.5  # J71: 92301.28999999999 -> 91888.5
$ws.Cells.Item(71, 12).Value2 = 275665.5  # L71: 276903.87 -> 275665.5
$ws.Cells.Item(71, 14).Value2 = -283153.5  # N71: -284391.87 -> -283153.5
$ws.Cells.Item(86, 8).Value2 = 2587.375  # H86: 2399.9 -> 2587.375
$ws.Cells.Item(86, 9).Value2 = 2528.4285  # I86: 2333.2222 -> 2528.4285
$ws.Cells.Item(86, 11).Value2 = 2528.4285  # K86: 2333.2222 -> 2528.4285
$ws.Cells.Item(86, 13).Value2 = -1405.4285  # M86: -1210.2222 -> -1405.4285
$ws.Cells.Item(89, 8).Value2 = 2587.375  # H89: 2399.9 -> 2587.375
$ws.Cells.Item(89, 9).Value2 = 2528.4285  # I89: 2333.2222 -> 2528.4285
$ws.Cells.Item(89, 11).Value2 = 12642.1425  # K89: 11666.111 -> 12642.1425
$ws.Cells.Item(89, 13).Value2 = -7026.1425  # M89: -6050.111000000001 -> -7026.1425
$ws.Cells.Item(107, 8).Value2 = 385.76923  # H107: 397.2 -> 385.76923
$ws.Cells.Item(107, 9).Value2 = 309.09525  # I107: 319.55 -> 309.09525
$ws.Cells.Item(107, 11).Value2 = 309.09525  # K107: 319.55 -> 309.09525
$ws.Cells.Item(107, 13).Value2 = 1610.90475  # M107: 1600.45 -> 1610.90475
$ws.Cells.Item(138, 8).Value2 = 73500  # H138: 74250 -> 73500
$ws.Cells.Item(138, 10).Value2 = 100000  # J138: 83333.336 -> 100000
$ws.Cells.Item(138, 12).Value2 = 100000  # L138: 83333.336 -> 100000
$ws.Cells.Item(138, 14).Value2 = -110280  # N138: -93613.336 -> -110280
$ws.Cells.Item(140, 8).Value2 = 89333.336  # H140: 86945 -> 89333.336
$ws.Cells.Item(140, 10).Value2 = 110000  # J140: 99926.664 -> 110000
$ws.Cells.Item(140, 12).Value2 = 110000  # L140: 99926.664 -> 110000
$ws.Cells.Item(140, 14).Value2 = -120360  # N140: -110286.664 -> -120360
$ws.Cells.Item(141, 8).Value2 = 63299.332  # H141: 57844.6 -> 63299.332
$ws.Cells.Item(141, 10).Value2 = 63299.332  # J141: 57844.6 -> 63299.332
$ws.Cells.Item(141, 12).Value2 = 63299.332  # L141: 57844.6 -> 63299.332
$ws.Cells.Item(141, 14).Value2 = -73659.33199999999  # N141: -68204.60000000001 -> -73659.33199999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value2 = 1983.3334  # H18: 2750 -> 1983.3334
$ws.Cells.Item(18, 9).Value2 = 475  # I18: 500 -> 475
$ws.Cells.Item(18, 11).Value2 = 1425  # K18: 1500 -> 1425
$ws.Cells.Item(18, 13).Value2 = -1256  # M18: -1331 -> -1256
$ws.Cells.Item(97, 8).Value2 = 1504.6666  # H97: 1507 -> 1504.6666
$ws.Cells.Item(97, 9).Value2 = 303.33334  # I97: 405 -> 303.33334
$ws.Cells.Item(97, 10).Value2 = 2706  # J97: 2058 -> 2706
$ws.Cells.Item(97, 11).Value2 = 910.0000200000001  # K97: 1215 -> 910.0000200000001
$ws.Cells.Item(97, 12).Value2 = 8118  # L97: 6174 -> 8118
$ws.Cells.Item(97, 13).Value2 = -414.0000200000001  # M97: -719 -> -414.0000200000001
$ws.Cells.Item(97, 14).Value2 = -9110  # N97: -7166 -> -9110
$ws.Cells.Item(107, 8).Value2 = 515.8333  # H107: 564.5 -> 515.8333
$ws.Cells.Item(107, 9).Value2 = 374  # I107: 429.8 -> 374
$ws.Cells.Item(107, 10).Value2 = 657.6667  # J107: 699.2 -> 657.6667
$ws.Cells.Item(107, 11).Value2 = 1122  # K107: 1289.4 -> 1122
$ws.Cells.Item(107, 12).Value2 = 1973.0001  # L107: 2097.6 -> 1973.0001
$ws.Cells.Item(107, 13).Value2 = 798  # M107: 630.5999999999999 -> 798
$ws.Cells.Item(107, 14).Value2 = -5813.0001  # N107: -5937.6 -> -5813.0001
$ws.Cells.Item(116, 8).Value2 = 2900  # H116: 2800 -> 2900
$ws.Cells.Item(116, 10).Value2 = 3000  # J116: 0 -> 3000
$ws.Cells.Item(116, 12).Value2 = 9000  # L116: 0 -> 9000
$ws.Cells.Item(116, 14).Value2 = -15884  # N116: None -> -15884
$ws.Cells.Item(124, 8).Value2 = 4989.6665  # H124: 4990 -> 4989.6665
$ws.Cells.Item(124, 9).Value2 = 4989.6665  # I124: 4990 -> 4989.6665
$ws.Cells.Item(124, 11).Value2 = 14968.9995  # K124: 14970 -> 14968.9995
$ws.Cells.Item(124, 13).Value2 = -10058.9995  # M124: -10060 -> -10058.9995
$ws.Cells.Item(139, 8).Value2 = 4219  # H139: 4627.125 -> 4219
$ws.Cells.Item(139, 9).Value2 = 4219  # I139: 4627.125 -> 4219
$ws.Cells.Item(139, 11).Value2 = 12657  # K139: 13881.375 -> 12657
$ws.Cells.Item(139, 13).Value2 = -7517  # M139: -8741.375 -> -7517

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value2 = 10000.5  # H21: 0 -> 10000.5
$ws.Cells.Item(21, 10).Value2 = 10000.5  # J21: 0 -> 10000.5
$ws.Cells.Item(21, 12).Value2 = 10000.5  # L21: 0 -> 10000.5
$ws.Cells.Item(21, 14).Value2 = -10346.5  # N21: None -> -10346.5
$ws.Cells.Item(30, 8).Value2 = 10000.5  # H30: 0 -> 10000.5
$ws.Cells.Item(30, 10).Value2 = 10000.5  # J30: 0 -> 10000.5
$ws.Cells.Item(30, 12).Value2 = 10000.5  # L30: 0 -> 10000.5
$ws.Cells.Item(30, 14).Value2 = -10210.5  # N30: None -> -10210.5
$ws.Cells.Item(132, 8).Value2 = 188834.83  # H132: 62643.527 -> 188834.83
$ws.Cells.Item(132, 9).Value2 = 277002.75  # I132: 72514.31 -> 277002.75
$ws.Cells.Item(132, 10).Value2 = 12499  # J132: 9999.333000000001 -> 12499
$ws.Cells.Item(132, 11).Value2 = 831008.25  # K132: 217542.93 -> 831008.25
$ws.Cells.Item(132, 12).Value2 = 37497  # L132: 29997.999 -> 37497
$ws.Cells.Item(132, 13).Value2 = -828478.25  # M132: -215012.93 -> -828478.25
$ws.Cells.Item(132, 14).Value2 = -42557  # N132: -35057.999 -> -42557

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value2 = 1877.4667  # H55: 2059.1538 -> 1877.4667
$ws.Cells.Item(55, 9).Value2 = 1833.4  # I55: 2192.25 -> 1833.4
$ws.Cells.Item(55, 10).Value2 = 1899.5  # J55: 2000 -> 1899.5
$ws.Cells.Item(55, 11).Value2 = 1833.4  # K55: 2192.25 -> 1833.4
$ws.Cells.Item(55, 12).Value2 = 1899.5  # L55: 2000 -> 1899.5
$ws.Cells.Item(55, 13).Value2 = -1660.4  # M55: -2019.25 -> -1660.4
$ws.Cells.Item(55, 14).Value2 = -2245.5  # N55: -2346 -> -2245.5
$ws.Cells.Item(82, 8).Value2 = 4721.231  # H82: 4973 -> 4721.231
$ws.Cells.Item(82, 9).Value2 = 3773.5  # I82: 4188.2 -> 3773.5
$ws.Cells.Item(82, 11).Value2 = 3773.5  # K82: 4188.2 -> 3773.5
$ws.Cells.Item(82, 13).Value2 = -3412.5  # M82: -3827.2 -> -3412.5
$ws.Cells.Item(85, 8).Value2 = 4721.231  # H85: 4973 -> 4721.231
$ws.Cells.Item(85, 9).Value2 = 3773.5  # I85: 4188.2 -> 3773.5
$ws.Cells.Item(85, 11).Value2 = 3773.5  # K85: 4188.2 -> 3773.5
$ws.Cells.Item(85, 13).Value2 = -2525.5  # M85: -2940.2 -> -2525.5
$ws.Cells.Item(132, 8).Value2 = 7498.75  # H132: 6399 -> 7498.75
$ws.Cells.Item(132, 9).Value2 = 0  # I132: 2000 -> 0
$ws.Cells.Item(132, 11).Value2 = 0  # K132: 6000 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: -3470 -> (removed)
$ws.Cells.Item(136, 8).Value2 = 6499.5  # H136: 5332.3335 -> 6499.5
$ws.Cells.Item(136, 9).Value2 = 2999  # I136: 2998.5 -> 2999
$ws.Cells.Item(136, 11).Value2 = 8997  # K136: 8995.5 -> 8997
$ws.Cells.Item(136, 13).Value2 = -6447  # M136: -6445.5 -> -6447

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value2 = 1245  # H6: 872.5 -> 1245
$ws.Cells.Item(6, 9).Value2 = 990  # I6: 663.3333 -> 990
$ws.Cells.Item(6, 11).Value2 = 990  # K6: 663.3333 -> 990
$ws.Cells.Item(6, 13).Value2 = -875  # M6: -548.3333 -> -875
$ws.Cells.Item(27, 8).Value2 = 100000  # H27: 0 -> 100000
$ws.Cells.Item(27, 10).Value2 = 100000  # J27: 0 -> 100000
$ws.Cells.Item(27, 12).Value2 = 100000  # L27: 0 -> 100000
$ws.Cells.Item(27, 14).Value2 = -100138  # N27: None -> -100138
$ws.Cells.Item(58, 8).Value2 = 22559.8  # H58: 28996.8 -> 22559.8
$ws.Cells.Item(58, 9).Value2 = 0  # I58: 50085 -> 0
$ws.Cells.Item(58, 10).Value2 = 22559.8  # J58: 23724.75 -> 22559.8
$ws.Cells.Item(58, 11).Value2 = 0  # K58: 50085 -> 0
$ws.Cells.Item(58, 12).Value2 = 22559.8  # L58: 23724.75 -> 22559.8
$ws.Cells.Item(58, 13).ClearContents()  # M58: -49777 -> (removed)
$ws.Cells.Item(58, 14).Value2 = -23175.8  # N58: -24340.75 -> -23175.8
$ws.Cells.Item(68, 8).Value2 = 30999.5  # H68: 31999 -> 30999.5
$ws.Cells.Item(68, 9).Value2 = 30000  # I68: 0 -> 30000
$ws.Cells.Item(68, 11).Value2 = 30000  # K68: 0 -> 30000
$ws.Cells.Item(68, 13).Value2 = -29189  # M68: None -> -29189
$ws.Cells.Item(71, 8).Value2 = 30999.5  # H71: 31999 -> 30999.5
$ws.Cells.Item(71, 9).Value2 = 30000  # I71: 0 -> 30000
$ws.Cells.Item(71, 11).Value2 = 90000  # K71: 0 -> 90000
$ws.Cells.Item(71, 13).Value2 = -85944  # M71: None -> -85944
$ws.Cells.Item(107, 8).Value2 = 953.6667  # H107: 1118.8 -> 953.6667
$ws.Cells.Item(107, 9).Value2 = 885.375  # I107: 961 -> 885.375
$ws.Cells.Item(107, 10).Value2 = 1500  # J107: 1750 -> 1500
$ws.Cells.Item(107, 11).Value2 = 2656.125  # K107: 2883 -> 2656.125
$ws.Cells.Item(107, 12).Value2 = 4500  # L107: 5250 -> 4500
$ws.Cells.Item(107, 13).Value2 = -736.125  # M107: -963 -> -736.125
$ws.Cells.Item(107, 14).Value2 = -8340  # N107: -9090 -> -8340
$ws.Cells.Item(115, 8).Value2 = 80000  # H115: 0 -> 80000
$ws.Cells.Item(115, 10).Value2 = 80000  # J115: 0 -> 80000
$ws.Cells.Item(115, 12).Value2 = 80000  # L115: 0 -> 80000
$ws.Cells.Item(115, 14).Value2 = -83134  # N115: None -> -83134
